$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.347.11"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.804.53"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D5").Value = "227.42"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "36.07"
$ws.Range("E8").Value = "  +10.71%  "
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "0.0964"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "2.065.53"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "11.75"
$ws.Range("E13").Value = "  +6.62%  "
$ws.Range("D14").Value = "1.805.19"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("E16").Value = "  +5.50%  "
$ws.Range("D17").Value = "34.334.65"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "69.10"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "245.51"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "11.58"
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "171.82"
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("E26").Value = "  +9.80%  "
$ws.Range("D27").Value = "16.89"
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "4.03"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").Value = "3.85"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").Value = "1.391.17"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").Value = "0.674"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "2.46"
$ws.Range("E37").Value = "  -5.83%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "0.966"
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "1.23"
$ws.Range("E41").Value = "  +10.81%  "
$ws.Range("D42").Value = "82.62"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "13.44"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("D48").Value = "1.965.89"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "104.53"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  +0.27%  "
